# The "Förändrad" (Changed) date column C was bumped by one day
# (2023-09-08 -> 2023-09-09, i.e. serial 45177 -> 45178) for every
# data row (rows 2-199) on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 199; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45177) {
        $cell.Value2 = 45178
    }
}
